# Applies the "Added periodic & upfront related scenarios" edit:
#  - ProductLoanInput!B17 (repaymentstrategy) changes from "RBI (India)"
#    to "Overdue/Due Fee/Int,Principal"
#  - ProductLoanInput!B3 (description) changes from the long auto-generated
#    description to "53-a"
#  - Selection on the input sheet ends up on B17

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Update the repayment strategy value
$wsInput.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Update the description value
$wsInput.Range("B3").Value = "53-a"

# The old description text required extra row height (wrapped text); the
# new short value no longer needs it, so restore the row to its default
# auto-fit height.
$wsInput.Rows.Item(3).EntireRow.AutoFit()

# Move the active selection to B17, matching the saved view state
$wsInput.Activate()
$wsInput.Range("B17").Select()
